$wb = $excel.ActiveWorkbook
$qa = $wb.Worksheets.Item("QA")

# Add the new "Datos" sheet right after "QA"
$ws2 = $wb.Worksheets.Add($null, $qa)
$ws2.Name = "Datos"

# Populate the province lookup list used by the data validation below
$values = @( `
    "provincia_nombre", `
    "[DESCONOCIDO]", `
    "COCLÉ", `
    "COLÓN", `
    "CHIRIQUÍ", `
    "DARIÉN", `
    "HERRERA", `
    "LOS SANTOS", `
    "PANAMÁ", `
    "VERAGUAS", `
    "COMARCA INDÍGENA GUNA YALA", `
    "COMARCA INDÍGENA EMBERÁ-WOUNAAN", `
    "COMARCA INDÍGENA NGÄBE-BUGLÉ", `
    "BOCAS DEL TORO", `
    "PANAMÁ OESTE" `
)
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws2.Cells.Item($i + 1, 1).Value = $values[$i]
}
$ws2.Range("A1:A15").Select()

# Update the ENF codes on the QA sheet (109/110 -> 112/113)
$qa.Range("A2").Value = "ENF-777-112"
$qa.Range("A3").Value = "ENF-777-113"

# Add list-based data validation on QA!G2:G3 sourced from the Datos sheet
$qa.Range("G2:G3").Validation.Add(3, 1, 1, 'Datos!$A$1:$A$15')

# Move the selection on QA and make it the active sheet again
$qa.Range("A5").Select()
$qa.Activate()
